$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E to fit the new, longer citation text (stored OOXML width
# ends up 5/6 of a character wider than the nominal ColumnWidth because of
# Excel's column-width-to-pixel quantization, so back that off here to land
# on a stored width of exactly 72).
$ws.Columns("E").ColumnWidth = 71.1666666666667

# Insert a new row at position 6 (pushes the old rows 6-7 down to 7-8),
# creating space for the "Supervisor" line that now sits between the
# MSc project line (row 5) and the "Mejor desempeño..." line (row 7).
$ws.Rows(6).Insert()
$ws.Rows(6).RowHeight = 30

# --- Row 2: PhD - Psychology ---------------------------------------------
$ws.Range("E2").Value = 'Proyecto de investigación: \href{https://dspace.stir.ac.uk/handle/1893/21102}{\textbf{\textit{Contextual musicality: vocal modulation and its perception in human social interaction}}}'

# --- Row 5: MSc in Evolutionary Psychology --------------------------------
$ws.Range("E5").Value = 'Proyecto de investigación: \textbf{\textit{Variation of pitch and loudness range of human voice in response to intra- and inter-sexual stimuli}}'

# --- Row 6 (new): Supervisor line, moved out of row 5 ---------------------
$ws.Range("E6").Value = 'Supervisor: \href{https://www.scraigroberts.com/}{Prof. S. Craig Roberts}'

# --- Row 7 (was row 6): Mejor desempeño... unchanged text, new position ---
$ws.Range("E7").Value = 'Mejor desempeño general en la maestría'

# --- Row 8 (was row 7): Licenciatura en Pedagogía Musical -----------------
$ws.Range("E8").Value = 'Trabajo de Grado: 4.90/5.00 | \href{https://revistas.pedagogica.edu.co/index.php/revistafba/article/view/50}{\textbf{\textit{El origen no humano de la música}}}'

# Restore the view state: scrolled down a couple of rows, with E9 selected
# (the cell just below the new last row of data).
$w = $excel.ActiveWindow
$w.ScrollRow = 3
$w.ScrollColumn = 1
$ws.Range("E9").Select()

Write-Host "edit applied"
